$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "91.276.69"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "3.152.94"
$ws.Range("E3").Value = "  +0.91%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'236.91"
$ws.Range("E5").Value = "  +7.94%  "
$ws.Range("D6").Value = "'641.63"
$ws.Range("E6").Value = "  +2.80%  "
$ws.Range("D7").Value = "'1.08"
$ws.Range("E7").Value = "  +8.83%  "
$ws.Range("D8").Value = "'0.373"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "3.151.00"
$ws.Range("E10").Value = "  +0.99%  "
$ws.Range("D11").Value = "'0.726"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("D13").Value = "'36.48"
$ws.Range("E13").Value = "  +4.85%  "
$ws.Range("D15").Value = "'5.61"
$ws.Range("E15").Value = "  +3.03%  "
$ws.Range("D16").Value = "90.776.30"
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "3.736.04"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").Value = "3.145.56"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "'3.75"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").Value = "'14.46"
$ws.Range("E21").Value = "  +2.28%  "
$ws.Range("D22").Value = "'451.17"
$ws.Range("E22").Value = "  +3.10%  "
$ws.Range("D23").Value = "'5.65"
$ws.Range("E23").Value = "  +9.01%  "
$ws.Range("D24").Value = "'9.08"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").Value = "'5.81"
$ws.Range("E25").Value = "  -7.04%  "
$ws.Range("D26").Value = "'91.70"
$ws.Range("E26").Value = "  +5.38%  "
$ws.Range("D27").Value = "'12.51"
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").Value = "3.302.23"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'9.77"
$ws.Range("E30").Value = "  +6.47%  "
$ws.Range("D31").Value = "'0.161"
$ws.Range("E31").Value = "  -4.67%  "
$ws.Range("D32").Value = "'0.993"
$ws.Range("E32").Value = "  +10.56%  "
$ws.Range("D33").Value = "'0.202"
$ws.Range("E33").Value = "  +31.88%  "
$ws.Range("D34").Value = "'27.06"
$ws.Range("E34").Value = "  +13.65%  "
$ws.Range("D35").Value = "'3.85"
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.151"
$ws.Range("E36").Value = "  +2.89%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "'515.24"
$ws.Range("E37").Value = "  -2.61%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D38").Value = "'1.95"
$ws.Range("E38").Value = "  +4.33%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").Value = "'7.15"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("D40").Value = "'1.32"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("E41").Value = "  +5.66%  "
$ws.Range("D42").Value = "'22.23"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").Value = "'0.0855"
$ws.Range("E43").Value = "  -4.24%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'3.37"
$ws.Range("E45").Value = "  +44.76%  "
$ws.Range("D46").Value = "'1.96"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").Value = "'0.704"
$ws.Range("E47").Value = "  +11.33%  "
$ws.Range("D48").Value = "'151.54"
$ws.Range("E48").Value = "  +1.54%  "
$ws.Range("D49").Value = "'45.63"
$ws.Range("E49").Value = "  +3.38%  "
$ws.Range("E50").Value = "  +7.54%  "
$ws.Range("E51").Value = "  +3.08%  "
